# Apply crypto symbol-list update (prices/volumes/coin reshuffle) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.99"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.45%"
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.63%"
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.65%"
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07617"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.10%"
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.247"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.92%"
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.603"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-8.59%"
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.82%"
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09873"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-11.89%"
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1767"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.83%"
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09126"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.93%"
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04364"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.40%"
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.16%"
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001260"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.76%"
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C15").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04162"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.43%"
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005815"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.89%"
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("B17").ClearFormats()
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.371"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.28%"
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("B18").ClearFormats()
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C18").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.437"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-5.83%"
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("B19").ClearFormats()
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C19").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3299"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.69%"
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "MCDex"
$ws.Range("B20").ClearFormats()
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C20").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.758"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-8.63%"
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("B21").ClearFormats()
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("C21").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1349"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.35%"
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("B22").ClearFormats()
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("C22").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2843"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "11.53%"
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001214"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.31%"
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004065"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-5.24%"
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.36%"
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003008"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.83%"
$ws.Range("E26").ClearFormats()

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02419"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-6.52%"
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.40%"
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007831"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.58%"
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1304"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.41%"
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007068"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.93%"
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001947"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.80%"
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008368"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.36%"
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3332"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.42%"
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006383"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.20%"
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.17%"
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-26.85%"
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.006335"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "87.74%"
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.17%"
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.17%"
$ws.Range("E51").ClearFormats()
